$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 10332.667
$ws.Range("J69").Value = 10332.667
$ws.Range("L69").Value = 30998.001
$ws.Range("N69").Value = -32746.001
$ws.Range("H72").Value = 10332.667
$ws.Range("J72").Value = 10332.667
$ws.Range("L72").Value = 92994.003
$ws.Range("N72").Value = -101730.003
$ws.Range("H116").Value = 7350.1113
$ws.Range("I116").Value = 6190.8335
$ws.Range("K116").Value = 6190.8335
$ws.Range("M116").Value = -2748.8335
$ws.Range("H132").Value = 1206.9269
$ws.Range("I132").Value = 665.3946999999999
$ws.Range("K132").Value = 1996.1841
$ws.Range("M132").Value = 533.8159000000001
$ws.Range("H138").Value = 2751.721
$ws.Range("J138").Value = 3021.6667
$ws.Range("L138").Value = 9065.000100000001
$ws.Range("N138").Value = -19345.0001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8265895
$ws.Range("I2").Value = 10101527
$ws.Range("K2").Value = 10101527
$ws.Range("M2").Value = -10101414
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H32").Value = 41001.22
$ws.Range("I32").Value = 42484.137
$ws.Range("K32").Value = 42484.137
$ws.Range("M32").Value = -42197.137
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H45").Value = 2158
$ws.Range("I45").Value = 1806
$ws.Range("K45").Value = 1806
$ws.Range("M45").Value = -1429
$ws.Range("H74").Value = 1990.7778
$ws.Range("I74").Value = 1098.174
$ws.Range("J74").Value = 3570
$ws.Range("K74").Value = 1098.174
$ws.Range("L74").Value = 3570
$ws.Range("M74").Value = -224.174
$ws.Range("N74").Value = -5318
$ws.Range("H77").Value = 1990.7778
$ws.Range("I77").Value = 1098.174
$ws.Range("J77").Value = 3570
$ws.Range("K77").Value = 5490.87
$ws.Range("L77").Value = 17850
$ws.Range("M77").Value = -1122.87
$ws.Range("N77").Value = -26586
$ws.Range("H116").Value = 8265895
$ws.Range("I116").Value = 10101527
$ws.Range("K116").Value = 10101527
$ws.Range("M116").Value = -10099233
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8265895
$ws.Range("I3").Value = 10101527
$ws.Range("K3").Value = 10101527
$ws.Range("M3").Value = -10101413
$ws.Range("H86").Value = 51386.35
$ws.Range("I86").Value = 1368.9333
$ws.Range("J86").Value = 201438.6
$ws.Range("K86").Value = 1368.9333
$ws.Range("L86").Value = 201438.6
$ws.Range("M86").Value = -245.9332999999999
$ws.Range("N86").Value = -203684.6
$ws.Range("H89").Value = 51386.35
$ws.Range("I89").Value = 1368.9333
$ws.Range("J89").Value = 201438.6
$ws.Range("K89").Value = 6844.666499999999
$ws.Range("L89").Value = 1007193
$ws.Range("M89").Value = -1228.666499999999
$ws.Range("N89").Value = -1018425
$ws.Range("H94").Value = 1638.7142
$ws.Range("I94").Value = 1564.6428
$ws.Range("K94").Value = 1564.6428
$ws.Range("M94").Value = -1113.6428
$ws.Range("H99").Value = 3259
$ws.Range("I99").Value = 2431.6667
$ws.Range("K99").Value = 2431.6667
$ws.Range("M99").Value = -933.6667000000002
$ws.Range("H109").Value = 60000
$ws.Range("J109").Value = 60000
$ws.Range("L109").Value = 60000
$ws.Range("N109").Value = -62774
$ws.Range("H134").Value = 6377.147
$ws.Range("I134").Value = 4886.5835
$ws.Range("J134").Value = 9954.5
$ws.Range("K134").Value = 14659.7505
$ws.Range("L134").Value = 29863.5
$ws.Range("M134").Value = -12124.7505
$ws.Range("N134").Value = -34933.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("H58").Value = 2977.58
$ws.Range("I58").Value = 2304.8125
$ws.Range("J58").Value = 4173.6113
$ws.Range("K58").Value = 2304.8125
$ws.Range("L58").Value = 4173.6113
$ws.Range("M58").Value = -2101.8125
$ws.Range("N58").Value = -4579.6113
$ws.Range("H105").Value = 929
$ws.Range("I105").Value = 712.6
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 712.6
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = 1034.4
$ws.Range("N105").Value = -5505
$ws.Range("H107").Value = 414
$ws.Range("I107").Value = 358.85715
$ws.Range("K107").Value = 358.85715
$ws.Range("M107").Value = 1561.14285
$ws.Range("H136").Value = 2977.58
$ws.Range("I136").Value = 2304.8125
$ws.Range("J136").Value = 4173.6113
$ws.Range("K136").Value = 6914.4375
$ws.Range("L136").Value = 12520.8339
$ws.Range("M136").Value = -4364.4375
$ws.Range("N136").Value = -17620.8339
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2000
$ws.Range("J52").Value = 2000
$ws.Range("L52").Value = 6000
$ws.Range("N52").Value = -6532
$ws.Range("H117").Value = 1072.1111
$ws.Range("J117").Value = 1526
$ws.Range("L117").Value = 4578
$ws.Range("N117").Value = -11462
$ws.Range("H129").Value = 50001250
$ws.Range("J129").Value = 83335310
$ws.Range("L129").Value = 250005930
$ws.Range("N129").Value = -250015930
$ws.Range("H131").Value = 13895690
$ws.Range("J131").Value = 9795.8125
$ws.Range("L131").Value = 29387.4375
$ws.Range("N131").Value = -39467.4375
$ws.Range("H140").Value = 1446.4
$ws.Range("I140").Value = 873
$ws.Range("K140").Value = 2619
$ws.Range("M140").Value = 2561
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 6873.75
$ws.Range("J17").Value = 7355.7144
$ws.Range("L17").Value = 7355.7144
$ws.Range("N17").Value = -7691.7144
$ws.Range("H97").Value = 790.8889
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 1000
$ws.Range("N97").Value = -1992
$ws.Range("H113").Value = 7126660
$ws.Range("I113").Value = 1523.25
$ws.Range("J113").Value = 14251797
$ws.Range("K113").Value = 1523.25
$ws.Range("L113").Value = 14251797
$ws.Range("M113").Value = 646.75
$ws.Range("N113").Value = -14256137
$ws.Range("H122").Value = 4218.067
$ws.Range("I122").Value = 948.6
$ws.Range("J122").Value = 10757
$ws.Range("K122").Value = 2845.8
$ws.Range("L122").Value = 32271
$ws.Range("M122").Value = -395.8000000000002
$ws.Range("N122").Value = -37171
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3450.25
$ws.Range("I22").Value = 2643
$ws.Range("K22").Value = 2643
$ws.Range("M22").Value = -2348
$ws.Range("H27").Value = 3450.25
$ws.Range("I27").Value = 2643
$ws.Range("K27").Value = 2643
$ws.Range("M27").Value = -2536
$ws.Range("H31").Value = 3665.4
$ws.Range("I31").Value = 716.25
$ws.Range("J31").Value = 5631.5
$ws.Range("K31").Value = 716.25
$ws.Range("L31").Value = 5631.5
$ws.Range("M31").Value = -468.25
$ws.Range("N31").Value = -6127.5
$ws.Range("H46").Value = 5136.9688
$ws.Range("I46").Value = 1075.25
$ws.Range("J46").Value = 5717.2144
$ws.Range("K46").Value = 1075.25
$ws.Range("L46").Value = 5717.2144
$ws.Range("M46").Value = -887.25
$ws.Range("N46").Value = -6093.2144
$ws.Range("H55").Value = 166
$ws.Range("I55").Value = 107.875
$ws.Range("J55").Value = 212.5
$ws.Range("K55").Value = 107.875
$ws.Range("L55").Value = 212.5
$ws.Range("M55").Value = 65.125
$ws.Range("N55").Value = -558.5
$ws.Range("H93").Value = 763.46155
$ws.Range("I93").Value = 737.5
$ws.Range("J93").Value = 850
$ws.Range("K93").Value = 737.5
$ws.Range("L93").Value = 850
$ws.Range("M93").Value = 510.5
$ws.Range("N93").Value = -3346
$ws.Range("H100").Value = 10420338
$ws.Range("I100").Value = 50002790
$ws.Range("K100").Value = 50002790
$ws.Range("M100").Value = -50002249
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 575
$ws.Range("I7").Value = 575
$ws.Range("K7").Value = 575
$ws.Range("M7").Value = -462
$ws.Range("H81").Value = 12971.037
$ws.Range("I81").Value = 4851.5835
$ws.Range("J81").Value = 19466.6
$ws.Range("K81").Value = 9703.166999999999
$ws.Range("L81").Value = 38933.2
$ws.Range("M81").Value = -8642.166999999999
$ws.Range("N81").Value = -41055.2
$ws.Range("H84").Value = 12971.037
$ws.Range("I84").Value = 4851.5835
$ws.Range("J84").Value = 19466.6
$ws.Range("K84").Value = 48515.835
$ws.Range("L84").Value = 194666
$ws.Range("M84").Value = -43211.835
$ws.Range("N84").Value = -205274
$ws.Range("H107").Value = 3580.2
$ws.Range("I107").Value = 1114.7142
$ws.Range("K107").Value = 3344.1426
$ws.Range("M107").Value = -1424.1426
